$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3748.25
$ws.Range("I40").Value = 3350
$ws.Range("K40").Value = 3350
$ws.Range("M40").Value = -3175
$ws.Range("H96").Value = 767.625
$ws.Range("I96").Value = 286.5
$ws.Range("J96").Value = 1248.75
$ws.Range("K96").Value = 859.5
$ws.Range("L96").Value = 3746.25
$ws.Range("M96").Value = 513.5
$ws.Range("N96").Value = -6492.25
$ws.Range("H137").Value = 2425.0435
$ws.Range("I137").Value = 2337.158
$ws.Range("K137").Value = 7011.474
$ws.Range("M137").Value = -4461.474
$ws.Range("H138").Value = 4105.631
$ws.Range("I138").Value = 1210.0834
$ws.Range("J138").Value = 5263.85
$ws.Range("K138").Value = 3630.2502
$ws.Range("L138").Value = 15791.55
$ws.Range("M138").Value = 1509.7498
$ws.Range("N138").Value = -26071.55

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 37206.5
$ws.Range("J28").Value = 53942
$ws.Range("L28").Value = 53942
$ws.Range("N28").Value = -54326
$ws.Range("H32").Value = 822.85
$ws.Range("I32").Value = 822.5204
$ws.Range("J32").Value = 839
$ws.Range("K32").Value = 822.5204
$ws.Range("L32").Value = 839
$ws.Range("M32").Value = -535.5204
$ws.Range("N32").Value = -1413
$ws.Range("H99").Value = 37206.5
$ws.Range("J99").Value = 53942
$ws.Range("L99").Value = 53942
$ws.Range("N99").Value = -59932
$ws.Range("H110").Value = 41667770
$ws.Range("I110").Value = 1258.5714
$ws.Range("K110").Value = 1258.5714
$ws.Range("M110").Value = 786.4286
$ws.Range("H132").Value = 5375.184
$ws.Range("I132").Value = 1732.3158
$ws.Range("J132").Value = 9018.053
$ws.Range("K132").Value = 5196.9474
$ws.Range("L132").Value = 27054.159
$ws.Range("M132").Value = -2666.9474
$ws.Range("N132").Value = -32114.159

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 18488.715
$ws.Range("I26").Value = 12579.833
$ws.Range("K26").Value = 12579.833
$ws.Range("M26").Value = -12287.833
$ws.Range("H98").Value = 53942
$ws.Range("J98").Value = 53942
$ws.Range("L98").Value = 53942
$ws.Range("N98").Value = -59932
$ws.Range("H107").Value = 75003730
$ws.Range("I107").Value = 187500690
$ws.Range("K107").Value = 187500690
$ws.Range("M107").Value = -187498770
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").Value = ""
$ws.Range("H134").Value = 5637.8125
$ws.Range("I134").Value = 1864.3043
$ws.Range("J134").Value = 9109.440000000001
$ws.Range("K134").Value = 5592.9129
$ws.Range("L134").Value = 27328.32
$ws.Range("M134").Value = -3057.9129
$ws.Range("N134").Value = -32398.32

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5677.6387
$ws.Range("I31").Value = 2418.8
$ws.Range("K31").Value = 2418.8
$ws.Range("M31").Value = -2123.8
$ws.Range("H34").Value = 5677.6387
$ws.Range("I34").Value = 2418.8
$ws.Range("K34").Value = 2418.8
$ws.Range("M34").Value = -2216.8
$ws.Range("H58").Value = 7357009.5
$ws.Range("I58").Value = 10870949
$ws.Range("J58").Value = 9682.591
$ws.Range("K58").Value = 10870949
$ws.Range("L58").Value = 9682.591
$ws.Range("M58").Value = -10870746
$ws.Range("N58").Value = -10088.591
$ws.Range("H62").Value = 9994.5
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 9994.5
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 9994.5
$ws.Range("M62").Value = ""
$ws.Range("N62").Value = -11242.5
$ws.Range("H65").Value = 9994.5
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 9994.5
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 49972.5
$ws.Range("M65").Value = ""
$ws.Range("N65").Value = -56212.5
$ws.Range("H99").Value = 7734.6665
$ws.Range("I99").Value = 3799
$ws.Range("K99").Value = 3799
$ws.Range("M99").Value = -2301
$ws.Range("H126").Value = 7734.6665
$ws.Range("I126").Value = 3799
$ws.Range("K126").Value = 11397
$ws.Range("M126").Value = -8927
$ws.Range("H136").Value = 7357009.5
$ws.Range("I136").Value = 10870949
$ws.Range("J136").Value = 9682.591
$ws.Range("K136").Value = 32612847
$ws.Range("L136").Value = 29047.773
$ws.Range("M136").Value = -32610297
$ws.Range("N136").Value = -34147.773
$ws.Range("H141").Value = 80712.625
$ws.Range("J141").Value = 80712.625
$ws.Range("L141").Value = 80712.625
$ws.Range("N141").Value = -91072.625

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H128").Value = 495000.5
$ws.Range("I128").Value = 495000.5
$ws.Range("K128").Value = 1485001.5
$ws.Range("M128").Value = -1480021.5

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").Value = ""
$ws.Range("H102").Value = 3195.6875
$ws.Range("I102").Value = 3017.3076
$ws.Range("K102").Value = 3017.3076
$ws.Range("M102").Value = -1395.3076

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6712.8823
$ws.Range("I7").Value = 4291.25
$ws.Range("J7").Value = 8865.444
$ws.Range("K7").Value = 4291.25
$ws.Range("L7").Value = 8865.444
$ws.Range("M7").Value = -4179.25
$ws.Range("N7").Value = -9089.444
$ws.Range("H22").Value = 1353.875
$ws.Range("I22").Value = 627.8889
$ws.Range("J22").Value = 2287.2856
$ws.Range("K22").Value = 627.8889
$ws.Range("L22").Value = 2287.2856
$ws.Range("M22").Value = -332.8889
$ws.Range("N22").Value = -2877.2856
$ws.Range("H27").Value = 1353.875
$ws.Range("I27").Value = 627.8889
$ws.Range("J27").Value = 2287.2856
$ws.Range("K27").Value = 627.8889
$ws.Range("L27").Value = 2287.2856
$ws.Range("M27").Value = -520.8889
$ws.Range("N27").Value = -2501.2856
$ws.Range("H40").Value = 4965.92
$ws.Range("I40").Value = 2751.7856
$ws.Range("K40").Value = 2751.7856
$ws.Range("M40").Value = -2615.7856
$ws.Range("H46").Value = 22225360
$ws.Range("J46").Value = 22225360
$ws.Range("L46").Value = 22225360
$ws.Range("N46").Value = -22225736
$ws.Range("H126").Value = 6712.8823
$ws.Range("I126").Value = 4291.25
$ws.Range("J126").Value = 8865.444
$ws.Range("K126").Value = 12873.75
$ws.Range("L126").Value = 26596.332
$ws.Range("M126").Value = -10403.75
$ws.Range("N126").Value = -31536.332

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = ""
$ws.Range("N93").Value = ""
$ws.Range("H100").Value = 1412.3
$ws.Range("I100").Value = 984.75
$ws.Range("J100").Value = 1697.3334
$ws.Range("K100").Value = 1969.5
$ws.Range("L100").Value = 3394.6668
$ws.Range("M100").Value = -1428.5
$ws.Range("N100").Value = -4476.6668
$ws.Range("H136").Value = 23813298
$ws.Range("I136").Value = 37037830
$ws.Range("J136").Value = 9133.267
$ws.Range("K136").Value = 111113490
$ws.Range("L136").Value = 27399.801
$ws.Range("M136").Value = -111110940
$ws.Range("N136").Value = -32499.801
